$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Data correction: rows 30-54 in column C were mistakenly tagged "QSR";
#    correct them to match the rest of the sheet ("Cinema_Cap"). Rewriting
#    the whole range lets the unused "QSR" shared-string entry be dropped.
$ws.Range("C30:C54").Value = "Cinema_Cap"

# 2. Re-apply the AutoFilter over the full data range (it had been left at
#    the header-only range A1:AM1; the sheet's real extent is A1:AM54).
$ws.AutoFilterMode = $false
$null = $ws.Range("A1:AM54").AutoFilter()

# 3. The (hidden) _FilterDatabase defined name should track the same new
#    range as the AutoFilter.
$names = $wb.Names
$first = $names.Item(1)
$first.RefersTo = '=Cinema!$A$1:$AM$54'

# 4. A fresh _FilterDatabase_0_0 defined name is created (mirrors the
#    pattern already present for _FilterDatabase_0), still pointing at the
#    original header-only range.
$ws.Names.Add("_xlnm._FilterDatabase_0_0", '=Cinema!$A$1:$AM$1')

# 5. Scroll the frozen pane back up to the top (topLeftCell A18 -> A2) and
#    move the active selection from F43 to E1.
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$null = $ws.Range("E1").Select()
